$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row values (string contents stay the same where unchanged,
# but columns B/C/D headers are reshuffled)
$ws.Range("B1").Value = "Both"
$ws.Range("C1").Value = "Seulement Individuelles"
$ws.Range("D1").Value = "Seulement copropriétés"

# Add "N/A" to the shared strings table first so it gets the lower new index
$ws.Range("C5").Value = "N/A"
$ws.Range("D5").Value = "N/A"

# Fix row label A3 wording and restore trailing parenthesis (added after N/A
# so it lands on the next shared-string slot)
$ws.Range("A3").Value = "2 (surface + chambres)"

# Fill in new column C and D values, plus remaining column B values
$ws.Range("C2").Value = 76932.8
$ws.Range("D2").Value = 41292.7

$ws.Range("C3").Value = 75894.7
$ws.Range("D3").Value = 39144.9

$ws.Range("C4").Value = 75894.600000000006
$ws.Range("D4").Value = 39137.699999999997

# Apply alignment styles
$ws.Range("A1:D1").HorizontalAlignment = -4108  # xlCenter
$ws.Range("B2:D5").HorizontalAlignment = -4108  # xlCenter
$ws.Range("A2:A5").HorizontalAlignment = -4131  # xlLeft

# Update selection to match final state
$ws.Range("C8").Select()

$wb.Save()
